# Fill in the new "F" / count column (D and E) for rows 2-10, and update
# the active selection to match (D2 active cell, D2:E10 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D10").Value = "F"
$ws.Range("E2:E10").Value = 1

$ws.Range("D2:E10").Select()
